$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EditViewTest")

$ws.Range("A2").Value = "testViewTwo"
$ws.Range("B2").Value = "testViewTwoEdited"
